$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $val) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = "Normal"
}

$ws.Range("D2").Value = "65.812.62"
$ws.Range("E2").Value = "  +2.03%  "
$ws.Range("D3").Value = "3.491.36"
$ws.Range("E3").Value = "  +0.95%  "
Set-TextValue "D5" "581.24"
$ws.Range("E5").Value = "  +0.70%  "
Set-TextValue "D6" "161.17"
$ws.Range("E6").Value = "  +1.89%  "
$ws.Range("E7").Value = "  +0.02%  "
Set-TextValue "D8" "0.608"
$ws.Range("E8").Value = "  +9.43%  "
$ws.Range("D9").Value = "3.493.47"
$ws.Range("E9").Value = "  +1.00%  "
Set-TextValue "D10" "7.31"
$ws.Range("E10").Value = "  -3.09%  "
$ws.Range("E11").Value = "  +0.95%  "
Set-TextValue "D12" "0.449"
$ws.Range("E12").Value = "  +1.01%  "
$ws.Range("D13").Value = "4.085.57"
$ws.Range("E13").Value = "  +0.68%  "
$ws.Range("E14").Value = "  -0.87%  "
$ws.Range("E15").Value = "  -0.45%  "
Set-TextValue "D16" "28.86"
$ws.Range("E16").Value = "  +3.85%  "
$ws.Range("D17").Value = "65.776.42"
$ws.Range("E17").Value = "  +1.91%  "
$ws.Range("D18").Value = "3.479.99"
$ws.Range("E18").Value = "  +0.80%  "
Set-TextValue "D19" "6.50"
$ws.Range("E19").Value = "  +1.11%  "
Set-TextValue "D20" "14.36"
$ws.Range("E20").Value = "  +0.01%  "
Set-TextValue "D21" "391.58"
$ws.Range("E21").Value = "  -0.83%  "
Set-TextValue "D22" "8.28"
$ws.Range("E22").Value = "  -2.42%  "
Set-TextValue "D23" "0.554"
$ws.Range("E23").Value = "  +1.56%  "
Set-TextValue "D24" "73.56"
$ws.Range("E24").Value = "  +0.92%  "
$ws.Range("E25").Value = "  +0.26%  "
$ws.Range("E26").Value = "  +1.89%  "
Set-TextValue "D27" "9.81"
$ws.Range("E27").Value = "  +1.83%  "
$ws.Range("E28").Value = "  -0.38%  "
$ws.Range("E29").Value = "  -0.06%  "
Set-TextValue "D30" "6.41"
$ws.Range("E30").Value = "  +4.50%  "
Set-TextValue "D31" "1.45"
$ws.Range("E31").Value = "  +4.55%  "
$ws.Range("E32").Value = "  +1.39%  "
Set-TextValue "D33" "23.79"
$ws.Range("E33").Value = "  +0.10%  "
Set-TextValue "D34" "6.53"
$ws.Range("E34").Value = "  -2.29%  "
Set-TextValue "D35" "0.999"
$ws.Range("E35").Value = "  +0.07%  "
Set-TextValue "D36" "7.18"
$ws.Range("E36").Value = "  +2.02%  "
$ws.Range("E37").Value = "  +5.38%  "
Set-TextValue "D38" "163.08"
$ws.Range("E38").Value = "  +1.51%  "
$ws.Range("E39").Value = "  +4.85%  "
$ws.Range("D40").Value = "3.088.10"
$ws.Range("E40").Value = "  +5.11%  "
Set-TextValue "D41" "0.0776"
$ws.Range("E41").Value = "  -1.18%  "
Set-TextValue "D42" "27.21"
$ws.Range("E42").Value = "  -1.26%  "
Set-TextValue "D43" "0.0323"
$ws.Range("E43").Value = "  +0.22%  "
Set-TextValue "D44" "4.55"
$ws.Range("E44").Value = "  +2.61%  "
Set-TextValue "D45" "43.12"
$ws.Range("E45").Value = "  +2.10%  "
Set-TextValue "D46" "0.782"
$ws.Range("E46").Value = "  +0.69%  "
Set-TextValue "D47" "25.83"
$ws.Range("E47").Value = "  +7.05%  "
$ws.Range("E48").Value = "  +3.12%  "
$ws.Range("E49").Value = "  +1.26%  "
Set-TextValue "D50" "6.75"
$ws.Range("E50").Value = "  +3.09%  "

# Row 51: Stellar -> Bittensor (full row replacement)
$ws.Range("B51").Value = "Bittensor"
$ws.Range("C51").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
Set-TextValue "D51" "313.09"
$ws.Range("E51").Value = "  +5.70%  "
